$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo in the "Wen You..." speaker string (missing opening paren before
# the Mark Prell link) and swap the F3/F4 cell contents so that the shared
# string table ends up with George Davis before Wen You, matching the target.
$wenYouText = "[Wen You](https://dataifa.github.io/difa-project/comingsoon.html), [Nichole Szembrot](https://dataifa.github.io/difa-project/comingsoon.html), [Mark Prell](https://dataifa.github.io/difa-project/comingsoon.html), [Bruce Weinberg](https://dataifa.github.io/difa-project/comingsoon.html)"
$georgeDavisText = "[George Davis](https://dataifa.github.io/difa-project/george_davis.html), [Joe Cummins](https://dataifa.github.io/difa-project/comingsoon.html)"

$ws.Range("F3").Value = $georgeDavisText
$ws.Range("F4").Value = $wenYouText
$ws.Range("F3").Value = $wenYouText
$ws.Range("F4").Value = $georgeDavisText

# Update the sheet view: drop the scrolled topLeftCell and move the
# selection from D13 to F3.
$ws.Range("F3").Select() | Out-Null
